# Apply the "A6 Request List" content update.
# The sample-group rows (10-16) are updated so that:
#   - the two standard references for "1. Краве Сирене" swap order
#   - the extra "2. Кокоши яйца" row (old row 13, "БДС ХЕР7") is removed
#     and every following group shifts up by one row
#   - row 15/16 gain new indicator/standard values, plus a note in F15
#   - row 18's footer text is left untouched (same text, OOXML just
#     renumbers the shared-string ids when strings are reordered)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "1. Краве Сирене" rows: swap the two standard numbers ---
$ws.Range("D10").Value = "БДС 0001 АР"
$ws.Range("D11").Value = "БНС 1234 АМ"

# --- "2. Кокоши яйца" keeps only its first row (row 12, unchanged). ---
# Former row 13 ("Яйчен тест 2" / "БДС ХЕР7") becomes the
# "3. Печен фъстък" row that used to live on row 14.
$ws.Range("B13").Value = "3. Печен фъстък"
$ws.Range("C13").Value = "Биреност"
$ws.Range("D13").Value = "БДС 7410Ж"

# Former row 15 ("4. Свинско месо" / "Ешерихия коли" / "БДС 753691")
# moves up to row 14, now paired with "Киселинност".
$ws.Range("B14").Value = "4. Свинско месо"
$ws.Range("C14").Value = "Киселинност"
$ws.Range("D14").Value = "БДС 456АЕЕР"

# Row 15 becomes a second "Свинско месо" line: label cleared, new
# indicator/standard values, and a remark in column F.
$ws.Range("B15").Value = ""
$ws.Range("C15").Value = "Киселинност"
$ws.Range("D15").Value = "БДС 7894"
$ws.Range("F15").Value = "Забележка ..."

# Row 16 gets the old "Ешерихия коли" indicator with a new standard.
$ws.Range("C16").Value = "Ешерихия коли"
$ws.Range("D16").Value = "БДС 788А"
